# Add a new row (17) of data to the "Daily" schedule sheet, matching the
# existing pattern of daily entries: Date, Daily.Measurements, Visual.Inspection,
# Initials (no Notes entry for this day).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daily")

$ws.Range("A17").Value = 20221014
$ws.Range("B17").Value = "completed"
$ws.Range("C17").Value = "completed"
$ws.Range("D17").Value = "AH"

$ws.Range("D21").Select()
